# The source sheet tracks weekly (daily-log) price observations for
# "Plátano" at "Terminal Hortofrutícola Agro Chillán". This commit adds a
# new pair of observations (Pintón / Primera Pintón) at the top of the
# existing data block (row 341), pushing the rest of the rows (formerly
# 341-392) down by two rows to 343-394. The dimension grows from
# A1:T392 to A1:T394 automatically once the new rows are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the first data row of the block
# (row 341). Excel's Insert() shifts rows 341:392 down to 343:394 and
# copies the formatting (incl. the date style on column D) from the row
# above, which matches the source data exactly.
$ws.Rows("341:342").Insert()

# New row 341 - "Pintón"
$ws.Range("A341").Value = 7
$ws.Range("B341").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C341").Value = "Ñuble"
$ws.Range("D341").Value = 44505
$ws.Range("E341").Value = 16
$ws.Range("F341").Value = "Fruta"
$ws.Range("G341").Value = 100108
$ws.Range("H341").Value = "Tropicales y subtropicales"
$ws.Range("I341").Value = 100108006
$ws.Range("J341").Value = "Plátano"
$ws.Range("K341").Value = "Sin especificar"
$ws.Range("L341").Value = "Pintón"
$ws.Range("M341").Value = 160
$ws.Range("N341").Value = 15000
$ws.Range("O341").Value = 16000
$ws.Range("P341").Value = 15500
$ws.Range("Q341").Value = "$/caja 20 kilos"
$ws.Range("R341").Value = "Ecuador"
$ws.Range("S341").Value = 775
$ws.Range("T341").Value = 20

# New row 342 - "Primera Pintón"
$ws.Range("A342").Value = 7
$ws.Range("B342").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C342").Value = "Ñuble"
$ws.Range("D342").Value = 44505
$ws.Range("E342").Value = 16
$ws.Range("F342").Value = "Fruta"
$ws.Range("G342").Value = 100108
$ws.Range("H342").Value = "Tropicales y subtropicales"
$ws.Range("I342").Value = 100108006
$ws.Range("J342").Value = "Plátano"
$ws.Range("K342").Value = "Sin especificar"
$ws.Range("L342").Value = "Primera Pintón"
$ws.Range("M342").Value = 240
$ws.Range("N342").Value = 17000
$ws.Range("O342").Value = 18000
$ws.Range("P342").Value = 17500
$ws.Range("Q342").Value = "$/caja 20 kilos"
$ws.Range("R342").Value = "Ecuador"
$ws.Range("S342").Value = 875
$ws.Range("T342").Value = 20

Write-Output "Inserted rows 341-342; dimension should now be A1:T394"
